# (BUG_FIX) Camera follows cat even when cat is colliding
#
# Three changes to the "Add Logic to assets" outline:
#  1. "Platforms"                                         -> cyan highlight
#  2. "Spawn platforms as the camera enters its area"      -> cyan highlight
#  3. "Spawn platforms in random places" / "Cat" / "Collider":
#       - "Spawn platforms in random places" highlight yellow -> cyan
#       - the "_GoBack" bookmark moves from the end of "Collider"
#         to the end of "Spawn platforms in random places"
#
# Paragraph.Range already includes the trailing paragraph mark, but this
# host's Range.HighlightColorIndex setter only ever touches run rPr (never
# the paragraph-mark rPr stored in pPr/rPr), so highlighting a whole
# paragraph that way leaves the pilcrow's own run-properties un-highlighted.
# Use Range.InsertXML (confirmed supported by this host's error text) to
# replace each target paragraph's XML outright with the exact OOXML we want,
# which lets us set both the run rPr *and* the paragraph-mark rPr, and lets
# us relocate the bookmark precisely.

$d = $word.ActiveDocument

function New-WordOpenXmlFragment([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" ' +
        'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---- locate the paragraphs we need to touch -------------------------------

$platformsPara = $null
$spawnEntersPara = $null
$spawnRandomPara = $null
$catPara = $null
$colliderPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match "^Platforms`r") { $platformsPara = $p }
    elseif ($t -match "^Spawn platforms as the camera enters its area`r") { $spawnEntersPara = $p }
    elseif ($t -match "^Spawn platforms in random places`r") { $spawnRandomPara = $p }
    elseif ($t -match "^Cat`r") { $catPara = $p }
    elseif ($t -match "^Collider`r") { $colliderPara = $p }
}

# ---- 1. "Platforms" gains a cyan highlight --------------------------------

$platformsXml = '<w:p w:rsidR="00AB614A" w:rsidRPr="000A6620" w:rsidRDefault="00AB614A" w:rsidP="00AB614A">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="000A6620"><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>Platforms</w:t></w:r></w:p>'
$platformsPara.Range.InsertXML((New-WordOpenXmlFragment $platformsXml))

# ---- 2. "Spawn platforms as the camera enters its area" gains a cyan highlight

$spawnEntersXml = '<w:p w:rsidR="000A6620" w:rsidRPr="00082B46" w:rsidRDefault="000A6620" w:rsidP="00AB614A">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00082B46"><w:rPr><w:highlight w:val="cyan"/></w:rPr>' +
    '<w:t>Spawn platforms as the camera enters its area</w:t></w:r></w:p>'
$spawnEntersPara.Range.InsertXML((New-WordOpenXmlFragment $spawnEntersXml))

# ---- 3. re-target highlight + move the "_GoBack" bookmark -----------------
# Re-fetch paragraphs 14-16 as one contiguous range (must still be adjacent
# to each other) and replace all three in one shot so the bookmark can move
# from the "Collider" paragraph to the "Spawn platforms in random places"
# paragraph.

$blockRange = $d.Range($spawnRandomPara.Range.Start, $colliderPara.Range.End)

$blockXml =
    '<w:p w:rsidR="00082B46" w:rsidRPr="00D70721" w:rsidRDefault="00082B46" w:rsidP="000B7F23">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr>' +
        '<w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr>' +
        '<w:r w:rsidRPr="00D70721"><w:rPr><w:highlight w:val="cyan"/></w:rPr>' +
        '<w:t>Spawn platforms in random places</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
    '<w:p w:rsidR="00AB614A" w:rsidRDefault="00AB614A" w:rsidP="00AB614A">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r><w:t>Cat</w:t></w:r>' +
    '</w:p>' +
    '<w:p w:rsidR="00AB614A" w:rsidRDefault="00AB614A" w:rsidP="00AB614A">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>' +
        '<w:rPr><w:highlight w:val="cyan"/></w:rPr></w:pPr>' +
        '<w:r w:rsidRPr="000A6620"><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>Collider</w:t></w:r>' +
    '</w:p>'

$blockRange.InsertXML((New-WordOpenXmlFragment $blockXml))
